# Initializing footer buttons component
# Restructure the decision-tree sheet: shift the "VALIDP" branch one column
# to the right to make room for a sibling, and add the new "VALIDS"/"VALIDF"
# branches (rows 55-74).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rows 6-10: shift the "next_statut" leaf one level to the right and
# introduce the isbtnValider leaf.
# ---------------------------------------------------------------------
$ws.Range("E6").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("F10").ClearContents()

$ws.Range("F6").Value = "next_statut == null"
$ws.Range("F6").Style = "Normal"

$ws.Range("G7").Value = "isbtnRenvoyer"
$ws.Range("G7").Style = "Normal"
$ws.Range("G7").Font.Italic = $true

$ws.Range("G8").Value = "isbtnValider"
$ws.Range("G8").Style = "Normal"
$ws.Range("G8").Font.Italic = $true

$ws.Range("F9").Value = "next_statut != null"
$ws.Range("F9").Style = "Normal"

$ws.Range("G10").Value = "next_statut == null"
$ws.Range("G10").Style = "Normal"
$ws.Range("G10").Font.Italic = $true

# ---------------------------------------------------------------------
# Rows 21-26 and 41-46: these leaves become italic (matching the style
# already used by their sibling leaves elsewhere in the tree).
# ---------------------------------------------------------------------
$ws.Range("F21").Font.Italic = $true
$ws.Range("F22").Font.Italic = $true
$ws.Range("F23").Font.Italic = $true
$ws.Range("F25").Font.Italic = $true
$ws.Range("F26").Font.Italic = $true
$ws.Range("F41").Font.Italic = $true
$ws.Range("F42").Font.Italic = $true
$ws.Range("F44").Font.Italic = $true
$ws.Range("F45").Font.Italic = $true
$ws.Range("F46").Font.Italic = $true

# ---------------------------------------------------------------------
# Row 29: the condition label wraps onto two lines now.
# ---------------------------------------------------------------------
$ws.Range("D29").WrapText = $true
$ws.Rows.Item(29).RowHeight = 30

# ---------------------------------------------------------------------
# New branches: statut == VALIDS (rows 55-64) and statut == VALIDF
# (rows 65-74), mirroring the existing VALIDP subtree.
# ---------------------------------------------------------------------
$ws.Range("B55").Value = "statut == VALIDS"

$ws.Range("C56").Value = "! currentUserValidS"

$ws.Range("D57").Value = "hasPermission (VALIDF)"

$ws.Range("E58").Value = "isbtnValiderF"
$ws.Range("E58").Font.Italic = $true

$ws.Range("E59").Value = "isbtnPLusDactions"
$ws.Range("E59").Font.Italic = $true

$ws.Range("F60").Value = "next_statut == null"

$ws.Range("G61").Value = "isbtnRenvoyer"
$ws.Range("G61").Font.Italic = $true

$ws.Range("D62").Value = "! hasPermission (VALIDF)"

$ws.Range("E63").Value = "isbtnOk"
$ws.Range("E63").Font.Italic = $true

$ws.Range("E64").Value = "isbtnPlusDactions"
$ws.Range("E64").Font.Italic = $true

$ws.Range("C65").Value = "currentUserValidS"

$ws.Range("D66").Value = "no Superior Validated"

$ws.Range("E67").Value = "isbtnAnnulerValider"
$ws.Range("E67").Font.Italic = $true

$ws.Range("E68").Value = "isbtnPlusDactions"
$ws.Range("E68").Font.Italic = $true

$ws.Range("D69").Value = "superior validated"

$ws.Range("E70").Value = "isbtnOk"
$ws.Range("E70").Font.Italic = $true

$ws.Range("E71").Value = "isbtnOptionsAnnuler"
$ws.Range("E71").Font.Italic = $true

$ws.Range("B72").Value = "statut == VALIDF"

$ws.Range("E72").Value = "isbtnPlusDactions"
$ws.Range("E72").Font.Italic = $true

$ws.Range("C73").Value = "isbtnOk"

$ws.Range("C74").Value = "isbtnOptionsAnnuler"

# ---------------------------------------------------------------------
# Column widths (best-fit, mirrors Excel auto-sizing after the new,
# wider content was added) and the freshly used columns F/G.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 11
$ws.Columns.Item(2).ColumnWidth = 15.42578125
$ws.Columns.Item(3).ColumnWidth = 18.42578125
$ws.Columns.Item(4).ColumnWidth = 34.42578125
$ws.Columns.Item(5).ColumnWidth = 20.42578125
$ws.Columns.Item(6).ColumnWidth = 17.5703125
$ws.Columns.Item(7).ColumnWidth = 14

# ---------------------------------------------------------------------
# Selection / scroll position reset to the top of the sheet.
# ---------------------------------------------------------------------
$ws.Range("B2").Select()
